# Auto-generated Excel COM-interop script to apply scheduled price/profit updates
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1118.9333
$ws.Range("J19").Value = 1347.2106
$ws.Range("L19").Value = 1347.2106
$ws.Range("N19").Value = -1697.2106
$ws.Range("H62").Value = 34485450
$ws.Range("I62").Value = 58825284
$ws.Range("J62").Value = 4013.6667
$ws.Range("K62").Value = 58825284
$ws.Range("L62").Value = 4013.6667
$ws.Range("M62").Value = -58824660
$ws.Range("N62").Value = -5261.6667
$ws.Range("H65").Value = 34485450
$ws.Range("I65").Value = 58825284
$ws.Range("J65").Value = 4013.6667
$ws.Range("K65").Value = 294126420
$ws.Range("L65").Value = 20068.3335
$ws.Range("M65").Value = -294123300
$ws.Range("N65").Value = -26308.3335
$ws.Range("H113").Value = 4289.7617
$ws.Range("I113").Value = 4100.4546
$ws.Range("J113").Value = 4498
$ws.Range("K113").Value = 4100.4546
$ws.Range("L113").Value = 4498
$ws.Range("M113").Value = -846.4546
$ws.Range("N113").Value = -11006
$ws.Range("H116").Value = 100003390
$ws.Range("I116").Value = 4245
$ws.Range("J116").Value = 500000000
$ws.Range("K116").Value = 4245
$ws.Range("L116").Value = 500000000
$ws.Range("M116").Value = -803
$ws.Range("N116").Value = -500006884
$ws.Range("H132").Value = 5454.8076
$ws.Range("I132").Value = 1534.375
$ws.Range("J132").Value = 52500
$ws.Range("K132").Value = 4603.125
$ws.Range("L132").Value = 157500
$ws.Range("M132").Value = -2073.125
$ws.Range("N132").Value = -162560
$ws.Range("H135").Value = 33334060
$ws.Range("I135").Value = 594.75
$ws.Range("J135").Value = 166667920
$ws.Range("K135").Value = 5352.75
$ws.Range("L135").Value = 1500011280
$ws.Range("M135").Value = -2817.75
$ws.Range("N135").Value = -1500016350
$ws.Range("H138").Value = 3255.0193
$ws.Range("I138").Value = 2506.4211
$ws.Range("J138").Value = 5286.9287
$ws.Range("K138").Value = 7519.263300000001
$ws.Range("L138").Value = 15860.7861
$ws.Range("M138").Value = -2379.263300000001
$ws.Range("N138").Value = -26140.7861

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21537.717
$ws.Range("I32").Value = 9051.424000000001
$ws.Range("K32").Value = 9051.424000000001
$ws.Range("M32").Value = -8764.424000000001
$ws.Range("H61").Value = 1494.5186
$ws.Range("I61").Value = 1284.5
$ws.Range("J61").Value = 2094.5715
$ws.Range("K61").Value = 1284.5
$ws.Range("L61").Value = 2094.5715
$ws.Range("M61").Value = -1072.5
$ws.Range("N61").Value = -2518.5715
$ws.Range("H74").Value = 35364.465
$ws.Range("I74").Value = 59990.766
$ws.Range("J74").Value = 3160.8462
$ws.Range("K74").Value = 59990.766
$ws.Range("L74").Value = 3160.8462
$ws.Range("M74").Value = -59116.766
$ws.Range("N74").Value = -4908.8462
$ws.Range("H77").Value = 35364.465
$ws.Range("I77").Value = 59990.766
$ws.Range("J77").Value = 3160.8462
$ws.Range("K77").Value = 299953.83
$ws.Range("L77").Value = 15804.231
$ws.Range("M77").Value = -295585.83
$ws.Range("N77").Value = -24540.231
$ws.Range("H110").Value = 1309.7273
$ws.Range("I110").Value = 1122.75
$ws.Range("J110").Value = 1416.5714
$ws.Range("K110").Value = 1122.75
$ws.Range("L110").Value = 1416.5714
$ws.Range("M110").Value = 922.25
$ws.Range("N110").Value = -5506.5714
$ws.Range("H122").Value = 1353.6904
$ws.Range("I122").Value = 1107.3715
$ws.Range("K122").Value = 3322.1145
$ws.Range("M122").Value = -872.1144999999997
$ws.Range("H132").Value = 1862.7576
$ws.Range("I132").Value = 1258.2593
$ws.Range("J132").Value = 4583
$ws.Range("K132").Value = 3774.7779
$ws.Range("L132").Value = 13749
$ws.Range("M132").Value = -1244.7779
$ws.Range("N132").Value = -18809
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H136").Value = 1494.5186
$ws.Range("I136").Value = 1284.5
$ws.Range("J136").Value = 2094.5715
$ws.Range("K136").Value = 3853.5
$ws.Range("L136").Value = 6283.7145
$ws.Range("M136").Value = -1303.5
$ws.Range("N136").Value = -11383.7145

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 609882.5600000001
$ws.Range("I134").Value = 692767.5
$ws.Range("J134").Value = 8966.75
$ws.Range("K134").Value = 2078302.5
$ws.Range("L134").Value = 26900.25
$ws.Range("M134").Value = -2075767.5
$ws.Range("N134").Value = -31970.25

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2432.9207
$ws.Range("I31").Value = 1093.7715
$ws.Range("K31").Value = 1093.7715
$ws.Range("M31").Value = -798.7715000000001
$ws.Range("H34").Value = 2432.9207
$ws.Range("I34").Value = 1093.7715
$ws.Range("K34").Value = 1093.7715
$ws.Range("M34").Value = -891.7715000000001
$ws.Range("H99").Value = 3086
$ws.Range("I99").Value = 3308.182
$ws.Range("K99").Value = 3308.182
$ws.Range("M99").Value = -1810.182
$ws.Range("H105").Value = 623.1
$ws.Range("I105").Value = 653.875
$ws.Range("K105").Value = 653.875
$ws.Range("M105").Value = 1093.125
$ws.Range("H122").Value = 684
$ws.Range("I122").Value = 572.63635
$ws.Range("K122").Value = 1717.90905
$ws.Range("M122").Value = 732.09095
$ws.Range("H126").Value = 3086
$ws.Range("I126").Value = 3308.182
$ws.Range("K126").Value = 9924.545999999998
$ws.Range("M126").Value = -7454.545999999998
$ws.Range("H134").Value = 2337.15
$ws.Range("I134").Value = 2090.6072
$ws.Range("J134").Value = 2912.4167
$ws.Range("K134").Value = 6271.821599999999
$ws.Range("L134").Value = 8737.250100000001
$ws.Range("M134").Value = -3736.821599999999
$ws.Range("N134").Value = -13807.2501

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4252.7144
$ws.Range("I3").Value = 3193.8
$ws.Range("K3").Value = 9581.400000000001
$ws.Range("M3").Value = -9469.400000000001
$ws.Range("H5").Value = 1025.7307
$ws.Range("I5").Value = 696
$ws.Range("J5").Value = 1200.2941
$ws.Range("K5").Value = 2088
$ws.Range("L5").Value = 3600.8823
$ws.Range("M5").Value = -1976
$ws.Range("N5").Value = -3824.8823
$ws.Range("H68").Value = 999.9
$ws.Range("I68").Value = 875
$ws.Range("J68").Value = 1083.1666
$ws.Range("K68").Value = 2625
$ws.Range("L68").Value = 3249.4998
$ws.Range("M68").Value = -1814
$ws.Range("N68").Value = -4871.4998
$ws.Range("H71").Value = 999.9
$ws.Range("I71").Value = 875
$ws.Range("J71").Value = 1083.1666
$ws.Range("K71").Value = 7875
$ws.Range("L71").Value = 9748.499400000001
$ws.Range("M71").Value = -3819
$ws.Range("N71").Value = -17860.4994
$ws.Range("H113").Value = 1443439.9
$ws.Range("I113").Value = 1894388
$ws.Range("J113").Value = 406
$ws.Range("K113").Value = 5683164
$ws.Range("L113").Value = 1218
$ws.Range("M113").Value = -5680994
$ws.Range("N113").Value = -5558
$ws.Range("H135").Value = 1025.7307
$ws.Range("I135").Value = 696
$ws.Range("J135").Value = 1200.2941
$ws.Range("K135").Value = 6264
$ws.Range("L135").Value = 10802.6469
$ws.Range("M135").Value = -3729
$ws.Range("N135").Value = -15872.6469

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 25642180
$ws.Range("I126").Value = 30304176
$ws.Range("J126").Value = 1195
$ws.Range("K126").Value = 90912528
$ws.Range("L126").Value = 3585
$ws.Range("M126").Value = -90910058
$ws.Range("N126").Value = -8525

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3288.7778
$ws.Range("I122").Value = 2899.5
$ws.Range("K122").Value = 8698.5
$ws.Range("M122").Value = -6248.5
$ws.Range("H132").Value = 5273.048
$ws.Range("I132").Value = 4923.875
$ws.Range("K132").Value = 14771.625
$ws.Range("M132").Value = -12241.625

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2750
$ws.Range("I62").Value = 2500
$ws.Range("K62").Value = 2500
$ws.Range("M62").Value = -1876
$ws.Range("H65").Value = 2750
$ws.Range("I65").Value = 2500
$ws.Range("K65").Value = 12500
$ws.Range("M65").Value = -9380
$ws.Range("H81").Value = 2825.25
$ws.Range("I81").Value = 2767
$ws.Range("J81").Value = 3000
$ws.Range("K81").Value = 5534
$ws.Range("L81").Value = 6000
$ws.Range("M81").Value = -4473
$ws.Range("N81").Value = -8122
$ws.Range("H84").Value = 2825.25
$ws.Range("I84").Value = 2767
$ws.Range("J84").Value = 3000
$ws.Range("K84").Value = 27670
$ws.Range("L84").Value = 30000
$ws.Range("M84").Value = -22366
$ws.Range("N84").Value = -40608
$ws.Range("H100").Value = 1712.1818
$ws.Range("I100").Value = 560.25
$ws.Range("J100").Value = 2370.4285
$ws.Range("K100").Value = 1120.5
$ws.Range("L100").Value = 4740.857
$ws.Range("M100").Value = -579.5
$ws.Range("N100").Value = -5822.857
$ws.Range("H126").Value = 981.55554
$ws.Range("I126").Value = 976.2857
$ws.Range("J126").Value = 1000
$ws.Range("K126").Value = 2928.8571
$ws.Range("L126").Value = 3000
$ws.Range("M126").Value = -458.8571000000002
$ws.Range("N126").Value = -7940
$ws.Range("H132").Value = 2184.7817
$ws.Range("I132").Value = 1983.2954
$ws.Range("J132").Value = 2990.7273
$ws.Range("K132").Value = 5949.8862
$ws.Range("L132").Value = 8972.1819
$ws.Range("M132").Value = -3419.8862
$ws.Range("N132").Value = -14032.1819
$ws.Range("H136").Value = 3190.2827
$ws.Range("I136").Value = 643.2778
$ws.Range("J136").Value = 12359.5
$ws.Range("K136").Value = 1929.8334
$ws.Range("L136").Value = 37078.5
$ws.Range("M136").Value = 620.1666
$ws.Range("N136").Value = -42178.5
